$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.226.23"
$ws.Range("E2").Value = "  -1.88%  "
$ws.Range("D3").Value = "1.820.53"
$ws.Range("E3").Value = "  -2.11%  "
$ws.Range("E4").Value = "  -1.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.21"
$ws.Range("E5").Value = "  -2.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  -1.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4275"
$ws.Range("E7").Value = "  -2.20%  "
$ws.Range("E8").Value = "  -2.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.90"
$ws.Range("E9").Value = "  -2.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07228"
$ws.Range("E10").Value = "  -2.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8604"
$ws.Range("E11").Value = "  -2.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.95"
$ws.Range("E12").Value = "  -2.87%  "
$ws.Range("D13").Value = "1.844.97"
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.657"
$ws.Range("E14").Value = "  -1.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07127"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.297"
$ws.Range("E16").Value = "  -3.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.46"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008863"
$ws.Range("E19").Value = "  -2.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  -1.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.04"
$ws.Range("E21").Value = "  -2.87%  "
$ws.Range("D22").Value = "27.268.38"
$ws.Range("E22").Value = "  -1.67%  "
$ws.Range("E23").Value = "  -2.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.88"
$ws.Range("E24").Value = "  -2.52%  "
$ws.Range("D25").Value = "2.056.58"
$ws.Range("E25").Value = "  -1.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.001"
$ws.Range("E26").Value = "  -1.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.29"
$ws.Range("E27").Value = "  -2.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.36"
$ws.Range("E28").Value = "  -1.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.114"
$ws.Range("E29").Value = "  +6.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.222"
$ws.Range("E30").Value = "  -3.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "116.18"
$ws.Range("E31").Value = "  -4.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08891"
$ws.Range("E32").Value = "  -1.86%  "
$ws.Range("E33").Value = "  -1.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7565"
$ws.Range("E34").Value = "  -1.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.447"
$ws.Range("E35").Value = "  -2.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.820"
$ws.Range("E36").Value = "  -7.13%  "
$ws.Range("E37").Value = "  -1.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.113"
$ws.Range("E38").Value = "  -2.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01967"
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05261"
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.902"
$ws.Range("E41").Value = "  +0.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.106"
$ws.Range("E42").Value = "  +2.13%  "
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5030"
$ws.Range("E44").Value = "  -2.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.598"
$ws.Range("E45").Value = "  -1.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.54"
$ws.Range("E46").Value = "  -2.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "106.65"
$ws.Range("E47").Value = "  -3.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4698"
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.004"
$ws.Range("E49").Value = "  -1.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06389"
$ws.Range("E50").Value = "  -1.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.659"
$ws.Range("E51").Value = "  -3.41%  "